$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("A2").Value = "9dd9db79-3ded-4362-bef2-4cdf58126809.md"
$ov.Range("B2").Value = "e2e\9dd9db79-3ded-4362-bef2-4cdf58126809.md"
$ov.Range("B2").Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e48a66f16c0b48588616d6a8127b28b2c1f92033/e2e/c1695c74-d002-440d-990b-85ecfacc4b2a.md", "", "", "e2e\9dd9db79-3ded-4362-bef2-4cdf58126809.md")
$ov.Range("G2").Value = "2016-09-06 05:16:36"

# --- zh-cn sheet ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("A2").Value = "9dd9db79-3ded-4362-bef2-4cdf58126809.md"
$zh.Range("A2").Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e48a66f16c0b48588616d6a8127b28b2c1f92033/e2e/c1695c74-d002-440d-990b-85ecfacc4b2a.md", "", "", "9dd9db79-3ded-4362-bef2-4cdf58126809.md")
$zh.Range("G2").Value = "9dd9db79-3ded-4362-bef2-4cdf58126809.3eda396770c4d04c7b581558dcf3bddf15e2063b.zh-cn.xlf"
$zh.Range("H2").Value = "2016-09-06 05:16:30"
$zh.Range("I2").Hyperlinks.Delete()
$zh.Range("I2").Style = "Normal"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"
$zh.Columns.Item(9).ColumnWidth = 18.6506053379604
$zh.Columns.Item(10).ColumnWidth = 21.7054770333426

# --- de-de sheet ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("A2").Value = "9dd9db79-3ded-4362-bef2-4cdf58126809.md"
$de.Range("A2").Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e48a66f16c0b48588616d6a8127b28b2c1f92033/e2e/c1695c74-d002-440d-990b-85ecfacc4b2a.md", "", "", "9dd9db79-3ded-4362-bef2-4cdf58126809.md")
$de.Range("G2").Value = "9dd9db79-3ded-4362-bef2-4cdf58126809.3eda396770c4d04c7b581558dcf3bddf15e2063b.de-de.xlf"
$de.Range("H2").Value = "2016-09-06 05:16:36"
$de.Range("I2").Hyperlinks.Delete()
$de.Range("I2").Style = "Normal"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"
$de.Columns.Item(9).ColumnWidth = 18.6506053379604
$de.Columns.Item(10).ColumnWidth = 21.7054770333426
